$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bug")

# --- Update row 14 column D: In-Progress -> Done ---
$ws.Range("D14").Value = "Done"

# --- New header for column E (bold, like the rest of the header row) ---
$ws.Range("E1").Value = "owner"
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)

# --- Owner column values (E2:E20) ---
$owners = @("sunil","sunil","saswat","saswat","saswat","saswat","nitin","saswat","saswat","saswat","nitin","sunil","nitin","saswat","nitin","nitin","sunil","sunil","saswat")
for ($i = 0; $i -lt $owners.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 5).Value = $owners[$i]
}

# --- New row 21 data ---
$ws.Range("A21").Value = 20
$ws.Range("B21").Value = "Payment Email and SMS text implementation"
$ws.Range("B20").Copy()
$ws.Range("B21").PasteSpecial(-4122)
$ws.Range("C21").Value = "code"
$ws.Range("D21").Value = "Pending"
$ws.Range("E21").Value = "Aakash"
$ws.Range("F21").Value = "Email template and sms text is pending from SNPL end"

# --- Size the newly introduced columns (bestFit-style) ---
$ws.Columns.Item(4).ColumnWidth = 10.16666666
$ws.Columns.Item(5).ColumnWidth = 6.45

# --- Update selection to match target state ---
$ws.Range("B22").Select()
